# Build site at 2023-04-12 14:53:07 UTC
# LOB1010.xlsx update:
#  - Fix "Objetivos" (row 10) text (was wrongly duplicated Daisy Rafaela da Silva text)
#  - Insert a new row for "Docentes responsaveis" (row 13) holding the teacher name
#    that previously (incorrectly) only appeared elsewhere
#  - Fix "Programa resumido" (now row 14) text
#  - Fix "Programa" (now row 16) text
#  - Fix "Metodo" / "Criterio" / "Norma de recuperacao" (rows 19-21) to hold the
#    correct cascading values
#  - Fix "Bibliografia" (now row 22) text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix "Objetivos:" (row 10) body text ---
$ws.Range("B10").Value = "GERAL: Dar noções gerais de direito, despertando o sentimento de cidadania através das garantias fundamentais asseguradas pela Constituição.`nESPECÍFICO: Preparar o aluno para o mercado de trabalho com as noções mínimas necessárias de direito relacionadas à sua profissão de engenheiro."
$ws.Range("C10").Value = "GERAL: Dar noções gerais de direito, despertando o sentimento de cidadania através das garantias fundamentais asseguradas pela Constituição.`nESPECÍFICO: Preparar o aluno para o mercado de trabalho com as noções mínimas necessárias de direito relacionadas à sua profissão de engenheiro."

# --- Insert a new row at 13 for "Docentes responsáveis:" value (teacher name) ---
$ws.Rows.Item(13).Insert()
# Row-insert copies column A's bold label style (and a phantom empty cell) down from
# row 12 into A13; this row has no label in column A, so reset it back to blank/plain.
$ws.Range("A13").Style = "Normal"
$ws.Range("A13").Value = $null
# Row-insert also copies column A's bold label style into B13 by default; restore the
# plain "data" look (non-bold, top-aligned, wrapped) used by every other B-column cell
# by copying the formatting from an existing data cell in column B.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("C13").Value = "6376612 - Daisy Rafaela da Silva"

# --- Fix "Programa resumido:" (now row 14) body text ---
$ws.Range("B14").Value = "Noções gerais de direito. O sistema constitucional brasileiro. Noções de direito civil. Propriedade intelectual. Noções de direito comercial e comercial internacional. Noções de direito administrativo. Noções de direito do trabalho. Noções de direito tributário. Regulamentação profissional."
$ws.Range("C14").Value = "Noções gerais de direito. O sistema constitucional brasileiro. Noções de direito civil. Propriedade intelectual. Noções de direito comercial e comercial internacional. Noções de direito administrativo. Noções de direito do trabalho. Noções de direito tributário. Regulamentação profissional."

# --- Fix "Programa:" (now row 16) body text ---
$programa = "01 - NOÇÕES GERAIS DE DIREITO: Orientação da disciplina. Conceito de Direito. Ato jurídico e ordem jurídica. Os vários ramos do Direito. Conceito de legislação. Jurisprudência e doutrina.`n02 - O SISTEMA CONSTITUCIONAL BRASILEIRO: Federação. República. Regime representativo. As garantias individuais.`n03 - NOÇÕES DE DIREITO CIVIL: Pessoas e bens. Direito de família. Atos jurídicos. Contratos. Atos ilícitos`n04 - NOÇÕES DE DIREITO COMERCIAL: Atos do comércio. Sociedades comerciais. Título de crédito. O cheque. A letra de câmbio. A nota promissória e a duplicata. Propriedade Industrial. Inventos, marcas e patentes. Proteção do direito autoral. . Contratos Comerciais: práticas ilegais e abusivas; regras de contratos internacionais. Da proteção ao consumidor`n05 - NOÇÕES DE DIREITO ADMINISTRATIVO: Ato administrativo e fato administrativo. Serviço público e de utilidade pública. Licitação e contrato administrativo.`n06 - NOÇÕES DE DIREITO DO TRABALHO: Conceitos fundamentais. Relações entre empregador e empregado. Higiene e segurança do trabalho. Previdência social. Justiça do trabalho.`n07 - NOÇÕES DE DIREITO TRIBUTÁRIO: O sistema tributário nacional. Tributos. Impostos, taxas e contribuições. Preços e tarifas.`n08 - REGULAMENTAÇÃO PROFISSIONAL: A garantia constitucional do exercício da profissão. A lei nº 5.194/66. Os órgãos regulamentadores da profissão. O exercício profissional. Atribuições. As atividades técnicas e econômicas da Engenharia. Responsabilidades decorrentes do exercício da engenharia."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# --- Fix "Método:" (now row 19) body text ---
$ws.Range("B19").Value = "Provas"
$ws.Range("C19").Value = "Provas"

# --- Fix "Critério:" (now row 20) body text ---
$ws.Range("B20").Value = "NF= (P1+P2)/2"
$ws.Range("C20").Value = "NF= (P1+P2)/2"

# --- Fix "Norma de recuperação:" (now row 21) body text ---
$ws.Range("B21").Value = "Reestudo com trabalhos e prova"
$ws.Range("C21").Value = "Reestudo com trabalhos e prova"

# --- Fix "Bibliografia:" (now row 22) body text ---
$bibliografia = "01. BATALHA, Wilson de Souza Campos. Introdução ao Estudo do Direito: Os Fundamentos e a Visão Histórica. Rio de Janeiro : Forense, 1981`n02. LIMA, Hermes. Introdução à Ciência do Direito. Rio de Janeiro: Freitas Bastos, 1980. `n03. NADER, Paulo. Introdução ao Estudo do Direito. Rio de Janeiro: Forense, 1982. `n04. PAUPERIO, A. Machado. Introdução ao Estudo do Direito. Rio de Janeiro: s.c.p., 1981.`n05. ROQUE, Ana. Direito Comercial Internacional. Portugal: Âncora Editora, 2004. `n06. PINHO, Rui Rebello; NASCIMENTO, Amauri Mascaro. Instituições de Direito Público e Privado. São Paulo: Atlas, 1984.`n07. REQUIAO, Rubens. Curso de Direito Comercial. São Paulo : Saraiva, 1981.`n08. BALEEIRO, Aliomar. Direito Tributário Brasileiro. Rio de Janeiro : Forense, 1981. `n09. BASTOS, Celso Ribeiro. Curso de Direito Constitucional.  São Paulo : Saraiva, 1989.`n10. ROMEIRO, José Antonio Nunes. Sociedade por Cotas de Responsabilidades Limitada. Curitiba: Juruá, 1984 `n11. RUSSOMANO, Mozart Victor. Comentários à Consolidação das Leis do Trabalho. Rio De Janeiro: Forense, 1994."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

Write-Host "LOB1010 worksheet updated"
